# Phoenix-Mesa-Scottsdale GDP data refresh ("updating GDP and algo")
# Updates the existing FRED GDP observations (2001-2019) with revised
# figures and appends the newly-released 2020 observation as row 31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised GDP values for existing years (B12:B30) ---
$ws.Range("B12").Value = 124351.109
$ws.Range("B13").Value = 130829.924
$ws.Range("B14").Value = 140576.011
$ws.Range("B15").Value = 150742.34
$ws.Range("B16").Value = 166970.628
$ws.Range("B17").Value = 180974.00899999999
$ws.Range("B18").Value = 191100.20800000001
$ws.Range("B19").Value = 190521.764
$ws.Range("B20").Value = 176581.61199999999
$ws.Range("B21").Value = 178397.03599999999
$ws.Range("B22").Value = 185788.92499999999
$ws.Range("B23").Value = 195823.019
$ws.Range("B24").Value = 201280.704
$ws.Range("B25").Value = 209291.46599999999
$ws.Range("B26").Value = 219956.59700000001
$ws.Range("B27").Value = 230743.32
$ws.Range("B28").Value = 243103.394
$ws.Range("B29").Value = 259286.47500000001
$ws.Range("B30").Value = 276914.30499999999

# --- New row: 2020-01-01 observation ---
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 281004.837
$ws.Range("B31").NumberFormat = "0.000"

# --- Leave the sheet with columns A:B selected (row 16 was last active) ---
$ws.Range("A16").Select()
$ws.Range("A1:B1048576").Select()
